$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some of the new "Price" strings look like plain numbers (e.g. "1.006", "310.01",
# "0.09170"). If we just assign them with Range.Value, Excel helpfully parses them as
# numbers, which both changes the cell type (t="n" instead of t="inlineStr"/shared
# string) and can silently normalise the text (e.g. "0.09170" -> 0.0917, dropping the
# trailing zero). The source file stores every Price/Volume cell as plain text, so for
# those cells we flip NumberFormat to Text ("@") right before writing the value, then
# call ClearFormats() to drop the temporary style again (the original cells carry no
# explicit style/"s" attribute at all).
$textForceCells = @(
    "D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11",
    "D12", "D13", "D14", "D15", "D16", "D19", "D20", "D22",
    "D23", "D25", "D26", "D27", "D28", "D29", "D30", "D31",
    "D32", "D34", "D35", "D36", "D37", "D38", "D39", "D40",
    "D41", "D42", "D43", "D45", "D46", "D47", "D49", "D50",
    "D51"
)
foreach ($cellRef in $textForceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# --- Row-by-row cell updates (Price + Volume columns; rows 41/42 also swap Coin/Link) ---
$ws.Range("D2").Value = '24.009.47'
$ws.Range("E2").Value = '  -0.86%  '
$ws.Range("D3").Value = '1.653.65'
$ws.Range("E3").Value = '  +0.41%  '
$ws.Range("D4").Value = '1.006'
$ws.Range("E4").Value = '  +0.67%  '
$ws.Range("D5").Value = '310.01'
$ws.Range("E5").Value = '  +0.12%  '
$ws.Range("D6").Value = '1.008'
$ws.Range("E6").Value = '  +0.72%  '
$ws.Range("D7").Value = '0.3927'
$ws.Range("E7").Value = '  +0.31%  '
$ws.Range("D8").Value = '0.3883'
$ws.Range("E8").Value = '  +0.39%  '
$ws.Range("D9").Value = '51.36'
$ws.Range("E9").Value = '  +3.68%  '
$ws.Range("D10").Value = '1.360'
$ws.Range("E10").Value = '  -0.07%  '
$ws.Range("D11").Value = '1.012'
$ws.Range("E11").Value = '  +0.83%  '
$ws.Range("D12").Value = '0.08472'
$ws.Range("E12").Value = '  -1.77%  '
$ws.Range("D13").Value = '23.97'
$ws.Range("E13").Value = '  +1.68%  '
$ws.Range("D14").Value = '7.219'
$ws.Range("E14").Value = '  +1.54%  '
$ws.Range("D15").Value = '7.895'
$ws.Range("E15").Value = '  +5.52%  '
$ws.Range("D16").Value = '0.00001317'
$ws.Range("E16").Value = '  +1.86%  '
$ws.Range("D17").Value = '1.649.23'
$ws.Range("E17").Value = '  +0.44%  '
$ws.Range("E18").Value = '  -0.56%  '
$ws.Range("D19").Value = '0.07021'
$ws.Range("E19").Value = '  +1.37%  '
$ws.Range("D20").Value = '20.11'
$ws.Range("E20").Value = '  -1.37%  '
$ws.Range("E21").Value = '  +0.22%  '
$ws.Range("D22").Value = '1.008'
$ws.Range("E22").Value = '  +0.64%  '
$ws.Range("D23").Value = '13.65'
$ws.Range("E23").Value = '  +0.61%  '
$ws.Range("D24").Value = '23.886.14'
$ws.Range("E24").Value = '  -1.31%  '
$ws.Range("D25").Value = '2.516'
$ws.Range("E25").Value = '  +4.56%  '
$ws.Range("D26").Value = '3.063'
$ws.Range("E26").Value = '  +7.89%  '
$ws.Range("D27").Value = '22.26'
$ws.Range("E27").Value = '  -0.51%  '
$ws.Range("D28").Value = '154.88'
$ws.Range("E28").Value = '  -1.92%  '
$ws.Range("D29").Value = '140.14'
$ws.Range("E29").Value = '  -0.31%  '
$ws.Range("D30").Value = '5.328'
$ws.Range("E30").Value = '  -0.81%  '
$ws.Range("D31").Value = '7.889'
$ws.Range("E31").Value = '  -6.47%  '
$ws.Range("D32").Value = '2.525'
$ws.Range("E32").Value = '  +4.88%  '
$ws.Range("D33").Value = '1.832.34'
$ws.Range("E33").Value = '  +0.95%  '
$ws.Range("D34").Value = '1.039'
$ws.Range("E34").Value = '  +9.26%  '
$ws.Range("D35").Value = '0.03029'
$ws.Range("E35").Value = '  +4.28%  '
$ws.Range("D36").Value = '0.08134'
$ws.Range("E36").Value = '  +0.31%  '
$ws.Range("D37").Value = '6.719'
$ws.Range("E37").Value = '  -3.45%  '
$ws.Range("D38").Value = '10.93'
$ws.Range("E38").Value = '  +7.74%  '
$ws.Range("D39").Value = '0.2717'
$ws.Range("E39").Value = '  +1.04%  '
$ws.Range("D40").Value = '0.09170'
$ws.Range("E40").Value = '  -0.41%  '
$ws.Range("B41").Value = 'Aptos'
$ws.Range("C41").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D41").Value = '13.53'
$ws.Range("E41").Value = '  +3.89%  '
$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").Value = '0.7526'
$ws.Range("E42").Value = '  -0.43%  '
$ws.Range("D43").Value = '1.425'
$ws.Range("E43").Value = '  -2.37%  '
$ws.Range("E44").Value = '  +1.79%  '
$ws.Range("D45").Value = '0.6952'
$ws.Range("E45").Value = '  +0.68%  '
$ws.Range("D46").Value = '2.477'
$ws.Range("E46").Value = '  +0.73%  '
$ws.Range("D47").Value = '4.088'
$ws.Range("E47").Value = '  -0.22%  '
$ws.Range("E48").Value = '  +0.66%  '
$ws.Range("D49").Value = '0.08303'
$ws.Range("E49").Value = '  -0.93%  '
$ws.Range("D50").Value = '134.57'
$ws.Range("E50").Value = '  +0.78%  '
$ws.Range("D51").Value = '1.411'
$ws.Range("E51").Value = '  +7.24%  '

# Drop the temporary Text format again so the Price cells end up unstyled, matching the
# original workbook (ClearFormats removes the style index but leaves the just-written
# string value/type untouched).
foreach ($cellRef in $textForceCells) {
    $ws.Range($cellRef).ClearFormats()
}
